$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the CHP column (C) and the solar-thermal column (which slides into C
# once CHP is removed) -- those energy sources are gone from this table.
$ws.Columns.Item(3).Delete()
$ws.Columns.Item(3).Delete()

# Remaining columns are now: B = net1, C = pvt1. Add a new "heat_pump1"
# column in D, matching the header formatting used by the other headers.
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)
$ws.Range("D1").Value = "Q_from_heat_pump1"
$ws.Range("D2").Value = "Q_heat_pump1_demand1"
$ws.Range("D3").Value = "Q_heat_pump1_net1"

# Add a new row describing the heat pump's "Q_to" flows, matching the
# formatting used for the other row labels in column A.
$ws.Range("A3").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("A4").Value = "Q_to_heat_pump1"
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 0
$ws.Range("D4").Value = 0
